$d = $word.ActiveDocument

$replacements = @(
    @{old="480÷4="; new="446÷3="},
    @{old="655÷4="; new="919÷8="},
    @{old="123÷6="; new="916÷2="},
    @{old="496÷4="; new="269÷8="},
    @{old="348÷9="; new="787÷5="},
    @{old="850÷5="; new="558÷5="},
    @{old="867÷5="; new="322÷8="},
    @{old="722÷2="; new="740÷8="},
    @{old="448÷6="; new="229÷6="},
    @{old="963÷3="; new="248÷4="},
    @{old="958÷7="; new="156÷7="},
    @{old="330÷4="; new="325÷4="},
    @{old="344÷6="; new="929÷8="},
    @{old="776÷4="; new="865÷8="},
    @{old="630÷9="; new="748÷8="},
    @{old="911÷4="; new="991÷6="},
    @{old="264÷9="; new="791÷9="},
    @{old="832÷6="; new="527÷5="},
    @{old="464÷9="; new="767÷9="},
    @{old="723÷3="; new="505÷9="},
    @{old="243÷4="; new="846÷5="},
    @{old="630÷4="; new="964÷5="},
    @{old="138÷2="; new="626÷7="},
    @{old="670÷4="; new="767÷8="},
    @{old="232÷7="; new="456÷4="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
